$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 30-32 (Pseudotime_1, Pseudotime_2, Pseudotime_3), shifting the
# remaining rows (Brain_region: ... through the last Mouse_ID row) up by 3.
$ws.Range("A30:I32").EntireRow.Delete()
